# hw6 pt 2 big upload
# Add a new "Sum- overall total population" column (F) that sums the
# male (D) and female (E) totals already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (F1) ---------------------------------------------------------
$ws.Range("F1").Value = "Sum- overall total population"
$ws.Range("F1").Style = "Normal"

# --- Data rows (F2:F101) --------------------------------------------------
# F2 is entered as a standalone formula (first cell typed by hand), then
# the rest of the column is filled down in two stretches, mirroring how
# the original workbook ended up with separate shared-formula groups.
$ws.Range("F2").Formula = "=SUM(D2,E2)"
$ws.Range("F2").Style = "Normal"

$ws.Range("F3:F66").Formula = "=SUM(D3,E3)"
$ws.Range("F3:F66").Style = "Normal"

$ws.Range("F67:F101").Formula = "=SUM(D67,E67)"
$ws.Range("F67:F101").Style = "Normal"

# --- Column width ----------------------------------------------------------
# (26.3 "characters" is the input value that rounds to the closest
# achievable stored width to the target 27.1640625)
$ws.Columns("F").ColumnWidth = 26.3

# --- Selection: whole column F, active cell F1 -----------------------------
[void]$ws.Columns("F").Select()
